$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New translation rows (Language | Key | Translation), continuing from row 218.
$rows = @(
    @{ Row = 219; Key = "lab.vendor.name.label.tooltip";        Translation = "Při zadávání jména výrobce se prosím snažte držet přesného názvu, včetně velikosti písmen a cizích znaků." },
    @{ Row = 220; Key = "lab.vendor.name.label";                 Translation = "Jméno výrobce" },
    @{ Row = 221; Key = "lab.vendor.create.submit";              Translation = "Vytvořit" },
    @{ Row = 222; Key = "lab.atomizer.create.submit";            Translation = "Vytvořit" },
    @{ Row = 223; Key = "lab.vendor.name.label.required";        Translation = "Jméno výrobce je povinné" },
    @{ Row = 224; Key = "lab.atomizer.name.label.required";      Translation = "Jméno atomizéru je povinné" },
    @{ Row = 225; Key = "lab.build.name.label.required";         Translation = "Jméno buildu je povinné" },
    @{ Row = 226; Key = "lab.build.create.submit";               Translation = "Vytvořit" },
    @{ Row = 227; Key = "lab.vendor.create.success";             Translation = "Výrobce [{{data.name}}] byl úspěšně vytvořen." },
    @{ Row = 228; Key = "lab.atomizer.create.success";           Translation = "Atomizér [{{data.name}}] byl úspěšně vytvořen." }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy formatting from the last pre-existing data row so the new rows
    # pick up the same style (wrapped-text, s="1") as the rest of the table.
    $src = $ws.Range("A218:C218")
    $dst = $ws.Range("A" + $rowNum + ":C" + $rowNum)
    $src.Copy()
    $dst.PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = "cs"
    $ws.Cells.Item($rowNum, 2).Value = $r.Key
    $ws.Cells.Item($rowNum, 3).Value = $r.Translation
}

# First new row wraps to two lines in the original workbook.
$ws.Rows.Item(219).RowHeight = 30

# Match the saved selection/view state.
$ws.Range("B219").Select() | Out-Null
